$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.256409
$ws.Range("H2").Value = 6.769227
$ws.Range("I2").Value = 0.08335947696385336
$ws.Range("J2").Value = 0.08335947696385337
$ws.Range("Q2").Value = 0.250973603843
$ws.Range("R2").Value = 2.258762434587
$ws.Range("S2").Value = 0.08335947696385336
$ws.Range("T2").Value = 0.08335947696385337

# Row 3
$ws.Range("I3").Value = 0.05236352458402688
$ws.Range("J3").Value = 0.05236352458402689
$ws.Range("S3").Value = 0.05236352458402688
$ws.Range("T3").Value = 0.05236352458402689

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.140981333333333
$ws.Range("H4").Value = 6.422944
$ws.Range("I4").Value = 0.0790951836019268
$ws.Range("J4").Value = 0.07909518360192681
$ws.Range("Q4").Value = 0.2381349307626667
$ws.Range("R4").Value = 2.143214376864
$ws.Range("S4").Value = 0.0790951836019268
$ws.Range("T4").Value = 0.07909518360192681

# Row 5
$ws.Range("G5").Value = 20.12912933333333
$ws.Range("H5").Value = 60.387388
$ws.Range("I5").Value = 0.7436389825445764
$ws.Range("J5").Value = 0.7436389825445765
$ws.Range("Q5").Value = 2.238902668358667
$ws.Range("R5").Value = 20.150124015228
$ws.Range("S5").Value = 0.7436389825445764
$ws.Range("T5").Value = 0.7436389825445765

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.124498666666667
$ws.Range("H6").Value = 3.373496
$ws.Range("I6").Value = 0.04154283230561651
$ws.Range("J6").Value = 0.04154283230561651
$ws.Range("Q6").Value = 0.1250746131973333
$ws.Range("R6").Value = 1.125671518776
$ws.Range("S6").Value = 0.04154283230561651
$ws.Range("T6").Value = 0.04154283230561651
